# Add a new "Machine 32" record (row 33) to the master-machine_master sheet,
# matching the new row 31 already present, and update the window
# selection/scroll state to reflect where the author left the cursor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new data row -------------------------------------------------
# Existing data runs from row 2 (id 10001) through row 32 (id 10031), so the
# new record becomes row 33 with id 10032.
$newRow = 33

$ws.Cells.Item($newRow, 1).Value = 10032                      # id
$ws.Cells.Item($newRow, 2).Value = "Machine 32"                # name
$ws.Cells.Item($newRow, 3).Value = "F4-30-B9-D4-CD-6F"         # mac_address
$ws.Cells.Item($newRow, 4).Value = "FB5962911665"               # serial_num
$ws.Cells.Item($newRow, 5).Value = "192.168.0.358"              # ip_address
$ws.Cells.Item($newRow, 6).Value = 1001                         # mspec_id
$ws.Cells.Item($newRow, 7).Value = "eng"                        # lang_code
$ws.Cells.Item($newRow, 8).Value = $true                        # is_active
$ws.Cells.Item($newRow, 9).Value = "superadmin"                 # cr_by
$ws.Cells.Item($newRow, 10).Value = "now()"                     # cr_dtimes
$ws.Cells.Item($newRow, 11).Value = "now()"                     # eff_dtimes

# --- Reflect the cursor / scroll position left behind by the edit ----------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C28").Select()
